# Update the Mantel correlogram table (recalculated with Euclidean
# distances).  Table layout: col1 = Distance Class (m), col2 = N,
# col3 = Mantel r, col4 = p.  Row 1 is the header; data rows are 2..15.
#
# Cells are addressed directly via Table.Cell(row, col) and their text is
# replaced by assigning Range.Text -- this only touches the run inside the
# cell (keeping its formatting) instead of a document-wide Find/Replace,
# which this engine resolves against the whole story rather than the
# supplied range.

$d = $word.ActiveDocument

function Set-CellText($rowIndex, $colIndex, $newText) {
    $table = $word.ActiveDocument.Tables.Item(1)
    $cell = $table.Cell($rowIndex, $colIndex)
    $cell.Range.Text = $newText
}

# Row 2 (1,250)
Set-CellText 2 3 "0.000"
Set-CellText 2 4 "0.461"

# Row 3 (3,750)
Set-CellText 3 3 "0.002"
Set-CellText 3 4 "0.909"

# Row 4 (6,250)
Set-CellText 4 3 "-0.021"
Set-CellText 4 4 "0.647"

# Row 5 (8,750)
Set-CellText 5 3 "0.008"
Set-CellText 5 4 "1"

# Row 6 (11,250)
Set-CellText 6 3 "-0.018"
Set-CellText 6 4 "1"

# Row 7 (13,750)
Set-CellText 7 3 "-0.024"
Set-CellText 7 4 "0.995"

# Row 8 (16,250)
Set-CellText 8 3 "-0.045"
Set-CellText 8 4 "0.266"

# Row 9 (18,750)
Set-CellText 9 3 "0.025"
Set-CellText 9 4 "1"

# Row 10 (21,250)
Set-CellText 10 3 "-0.039"
Set-CellText 10 4 "0.315"

# Row 11 (23,750)
Set-CellText 11 3 "0.010"
Set-CellText 11 4 "1"

# Row 12 (26,250)
Set-CellText 12 3 "0.032"
Set-CellText 12 4 "0.89"

# Row 13 (28,750) -- the p-value cell was bold (significant at p<0.05);
# the recalculated p (0.919) is no longer significant, so un-bold it too.
Set-CellText 13 3 "0.032"
Set-CellText 13 4 "0.919"
$table = $word.ActiveDocument.Tables.Item(1)
$cell = $table.Cell(13, 4)
$cell.Range.Font.Bold = $false

# Row 14 (31,250)
Set-CellText 14 3 "0.032"
Set-CellText 14 4 "0.967"

# Row 15 (33,750) -- only the Mantel r value changes; p stays "1"
Set-CellText 15 3 "-0.033"
